# Add the path to the image for the "Permanent Residence Permit" document
# (column P = "image") on the data row (row 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = "C:Users/vano/Documents/GitHub/ZPI_VAF/iaff_assistant/images/Cards/permanent.jpg"

# Match the resulting selection/active cell as seen in the saved file.
$ws.Range("P2").Select()
